$d = $word.ActiveDocument

# Each call below targets one, unambiguous occurrence and replaces only that
# single match (wdReplaceOne = 1), so the four edits can never bleed into one
# another even though several share the substring "July 01, 2022".

# 1) "...appeared in Court for sentencing on July 01, 2022." -> "...July 02, 2022."
$d.Content.Find.Execute(" on July 01, 2022.", $false, $false, $false, $false, $false,
                         $true, 1, $false, " on July 02, 2022.", 1)

# 2) Bold "July 01, 2022" (the pay-fines-and-costs-by date) -> "July 02, 2022"
#    Searching the bare date (not spanning into the surrounding plain-text run)
#    keeps the match inside the single bold run, so bold formatting survives.
#    By this point the only earlier "July 01, 2022" has already been replaced in
#    step 1, so wdReplaceOne (1) is guaranteed to hit this bold occurrence next.
$d.Content.Find.Execute("July 01, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "July 02, 2022", 1)

# 3) "August 30, 2022" -> "August 31, 2022"
$d.Content.Find.Execute("August 30, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "August 31, 2022", 1)

# 4) " license is suspended from July 01, 2022" -> " license is suspended from July 02, 2022"
$d.Content.Find.Execute(" license is suspended from July 01, 2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, " license is suspended from July 02, 2022", 1)
